# Changes 2 to 17
#
# - Cell A3 on Sheet1 changes from 2 to 17 (A6 = SUM(A2:A5) recalculates
#   automatically from 10 to 25).
# - The selected cell/active cell on Sheet1 moves from C6 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 17

$ws.Range("B3").Select() | Out-Null
